$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 7 new "Assignment / Grade / Comments" groups after the existing ones
# (columns AD:BD), each preceded by a spacer column, mirroring the layout
# already used in columns B:AB. Also raise the grade for the "7 CPP"
# assignment (now AD:AF) and add its comment.
# ---------------------------------------------------------------------------

$groups = @(
    @{ Assign = "AD"; Grade = "AE"; Comments = "AF"; AssignVal = "7 CPP";  GradeVal = 90;  CommentVal = "Very good!" },
    @{ Assign = "AH"; Grade = "AI"; Comments = "AJ"; AssignVal = "8 CPP";  GradeVal = 100; CommentVal = "Excellent!" },
    @{ Assign = "AL"; Grade = "AM"; Comments = "AN"; AssignVal = "9 CPP";  GradeVal = 100; CommentVal = "Excellent!" },
    @{ Assign = "AP"; Grade = "AQ"; Comments = "AR"; AssignVal = "10 CPP"; GradeVal = 90;  CommentVal = "Very good! (didn’t implement copy assignment operator properly - not copying values)" },
    @{ Assign = "AT"; Grade = "AU"; Comments = "AV"; AssignVal = "11 CPP"; GradeVal = 90;  CommentVal = "Very good (read my comments)" },
    @{ Assign = "AX"; Grade = "AY"; Comments = "AZ"; AssignVal = "12 CPP"; GradeVal = 100; CommentVal = "Excellent!" },
    @{ Assign = "BB"; Grade = "BC"; Comments = "BD"; AssignVal = "Final Project"; GradeVal = 100; CommentVal = "Excellent!" }
)

$spacers = @("AC", "AG", "AK", "AO", "AS", "AW", "BA")

# ---- Header row (row 1): Assignment / Grade / Comments labels ----
foreach ($g in $groups) {
    $ws.Range($g.Assign + "1").Value = "Assignment"
    $ws.Range($g.Grade + "1").Value = "Grade"
    $ws.Range($g.Comments + "1").Value = "Comments"

    $ws.Range("Z1:AB1").Copy()
    $ws.Range($g.Assign + "1:" + $g.Comments + "1").PasteSpecial(-4122)
}

foreach ($sp in $spacers) {
    $ws.Range($sp + "1").Value = $null
    $ws.Range("Y1").Copy()
    $ws.Range($sp + "1").PasteSpecial(-4122)
}

# ---- Data row (row 2): values for each new group ----
foreach ($g in $groups) {
    $ws.Range($g.Assign + "2").Value = $g.AssignVal
    $ws.Range($g.Grade + "2").Value = $g.GradeVal
    $ws.Range($g.Comments + "2").Value = $g.CommentVal

    $ws.Range("Z2:AB2").Copy()
    $ws.Range($g.Assign + "2:" + $g.Comments + "2").PasteSpecial(-4122)
}

foreach ($sp in $spacers) {
    $ws.Range($sp + "2").Value = $null
    $ws.Range("Y1").Copy()
    $ws.Range($sp + "2").PasteSpecial(-4122)
}

# ---- Blank rows 3-6: extend borders/styles into the new columns ----
foreach ($r in 3..6) {
    foreach ($g in $groups) {
        $ws.Range("Z" + $r + ":AB" + $r).Copy()
        $ws.Range($g.Assign + $r + ":" + $g.Comments + $r).PasteSpecial(-4122)
    }
    foreach ($sp in $spacers) {
        $ws.Range("Y1").Copy()
        $ws.Range($sp + $r).PasteSpecial(-4122)
    }
}

# ---------------------------------------------------------------------------
# Row heights: header row now wraps onto two lines, and the data row grew
# tall enough to show the longest multi-line comment in full.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(2).RowHeight = 180

# ---------------------------------------------------------------------------
# Selection / view: the sheet was scrolled right and BA1:BA6 selected.
# ---------------------------------------------------------------------------
$ws.Range("BA1:BA6").Select()
